$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Handoff/Handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-17 04:58:54"
$wsZhCn.Range("G2").Value = "2016-02-17 04:59:40"
$wsZhCn.Range("D3").Value = "2016-02-17 04:58:54"
$wsZhCn.Range("G3").Value = "2016-02-17 04:59:40"

# de-de sheet: update Handoff/Handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-17 04:59:04"
$wsDeDe.Range("G2").Value = "2016-02-17 04:59:57"
$wsDeDe.Range("D3").Value = "2016-02-17 04:59:04"
$wsDeDe.Range("G3").Value = "2016-02-17 04:59:57"
